$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.558.78'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.934.48'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.62'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4842'
$ws.Range('E7').Value = '  +2.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2915'
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06794'
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '113.27'
$ws.Range('E10').Value = '  +6.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.47'
$ws.Range('E11').Value = '  +4.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.936.30'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.506'
$ws.Range('E13').Value = '  +2.86%  '
$ws.Range('E14').Value = '  -1.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6794'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '299.33'
$ws.Range('E16').Value = '  +3.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.585.32'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.09'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007646'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9992'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.188.61'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.553'
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9997'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.519'
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.573'
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.38'
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.42'
$ws.Range('E27').Value = '  -1.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.128'
$ws.Range('E28').Value = '  -0.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1069'
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.432'
$ws.Range('E30').Value = '  +1.73%  '
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.094'
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04999'
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7511'
$ws.Range('E34').Value = '  +1.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.147'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('E37').Value = '  -1.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.692'
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.032'
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '110.16'
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4462'
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8723'
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.828'
$ws.Range('E43').Value = '  -1.58%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '69.77'
$ws.Range('E45').Value = '  +2.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.312'
$ws.Range('E46').Value = '  -0.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '49.30'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.308'
$ws.Range('E48').Value = '  -1.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1233'
$ws.Range('E49').Value = '  -2.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.2549'
$ws.Range('E50').Value = '  +2.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.10'
$ws.Range('E51').Value = '  -0.74%  '
